$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 300.4762
$ws.Range("I5").Value = 120.36842
$ws.Range("J5").Value = 2011.5
$ws.Range("K5").Value = 120.36842
$ws.Range("L5").Value = 2011.5
$ws.Range("M5").Value = -5.36842
$ws.Range("N5").Value = -2241.5
$ws.Range("H38").Value = 2187.889
$ws.Range("J38").Value = 2900
$ws.Range("L38").Value = 8700
$ws.Range("N38").Value = -9444
$ws.Range("H111").Value = 90045.25
$ws.Range("I111").Value = 149079.28
$ws.Range("K111").Value = 447237.84
$ws.Range("M111").Value = -444170.84
$ws.Range("H116").Value = 13513.25
$ws.Range("I116").Value = 4650.8
$ws.Range("K116").Value = 4650.8
$ws.Range("M116").Value = -1208.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1532.9333
$ws.Range("I2").Value = 1545.6923
$ws.Range("K2").Value = 1545.6923
$ws.Range("M2").Value = -1432.6923
$ws.Range("H32").Value = 4059.975
$ws.Range("I32").Value = 4059.975
$ws.Range("K32").Value = 4059.975
$ws.Range("M32").Value = -3772.975
$ws.Range("H110").Value = 239850.19
$ws.Range("I110").Value = 279657.84
$ws.Range("K110").Value = 279657.84
$ws.Range("M110").Value = -277612.84
$ws.Range("H116").Value = 1532.9333
$ws.Range("I116").Value = 1545.6923
$ws.Range("K116").Value = 1545.6923
$ws.Range("M116").Value = 748.3077000000001
$ws.Range("H122").Value = 6534.091
$ws.Range("I122").Value = 8050.125
$ws.Range("J122").Value = 5667.7856
$ws.Range("K122").Value = 24150.375
$ws.Range("L122").Value = 17003.3568
$ws.Range("M122").Value = -21700.375
$ws.Range("N122").Value = -21903.3568

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1532.9333
$ws.Range("I3").Value = 1545.6923
$ws.Range("K3").Value = 1545.6923
$ws.Range("M3").Value = -1431.6923

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3858.75
$ws.Range("I62").Value = 2478.3333
$ws.Range("K62").Value = 2478.3333
$ws.Range("M62").Value = -1854.3333
$ws.Range("H65").Value = 3858.75
$ws.Range("I65").Value = 2478.3333
$ws.Range("K65").Value = 12391.6665
$ws.Range("M65").Value = -9271.666499999999
$ws.Range("H99").Value = 560617.1
$ws.Range("I99").Value = 5311
$ws.Range("J99").Value = 1254749.8
$ws.Range("K99").Value = 5311
$ws.Range("L99").Value = 1254749.8
$ws.Range("M99").Value = -3813
$ws.Range("N99").Value = -1257745.8
$ws.Range("H126").Value = 560617.1
$ws.Range("I126").Value = 5311
$ws.Range("J126").Value = 1254749.8
$ws.Range("K126").Value = 15933
$ws.Range("L126").Value = 3764249.4
$ws.Range("M126").Value = -13463
$ws.Range("N126").Value = -3769189.4
$ws.Range("H132").Value = 1118.8572
$ws.Range("I132").Value = 888.6667
$ws.Range("K132").Value = 2666.0001
$ws.Range("M132").Value = -136.0001000000002
$ws.Range("H134").Value = 324974.3
$ws.Range("I134").Value = 2534.6206
$ws.Range("K134").Value = 7603.861800000001
$ws.Range("M134").Value = -5068.861800000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 124.166664
$ws.Range("I2").Value = 43.666668
$ws.Range("J2").Value = 164.41667
$ws.Range("K2").Value = 262.000008
$ws.Range("L2").Value = 986.5000200000001
$ws.Range("M2").Value = -149.000008
$ws.Range("N2").Value = -1212.50002
$ws.Range("H33").Value = 3458040.5
$ws.Range("J33").Value = 125110.625
$ws.Range("L33").Value = 750663.75
$ws.Range("N33").Value = -751229.75
$ws.Range("H92").Value = 417441.66
$ws.Range("I92").Value = 769728.7
$ws.Range("J92").Value = 1102.4546
$ws.Range("K92").Value = 2309186.1
$ws.Range("L92").Value = 3307.3638
$ws.Range("M92").Value = -2307938.1
$ws.Range("N92").Value = -5803.3638
$ws.Range("H93").Value = 6999
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 6999
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 20997
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -24741
$ws.Range("H95").Value = 20000
$ws.Range("I95").Value = 20000
$ws.Range("K95").Value = 60000
$ws.Range("M95").Value = -57941
$ws.Range("H96").Value = 338666.66
$ws.Range("J96").Value = 502000
$ws.Range("L96").Value = 1506000
$ws.Range("N96").Value = -1510118
$ws.Range("H99").Value = 2663.1667
$ws.Range("I99").Value = 1489.5
$ws.Range("J99").Value = 3250
$ws.Range("K99").Value = 4468.5
$ws.Range("L99").Value = 9750
$ws.Range("M99").Value = -2222.5
$ws.Range("N99").Value = -14242
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("N104").ClearContents()
$ws.Range("H107").Value = 79640.03999999999
$ws.Range("I107").Value = 1005.5455
$ws.Range("J107").Value = 137305.33
$ws.Range("K107").Value = 3016.6365
$ws.Range("L107").Value = 411915.99
$ws.Range("M107").Value = -1096.6365
$ws.Range("N107").Value = -415755.99
$ws.Range("H132").Value = 1574496.1
$ws.Range("J132").Value = 2504250.8
$ws.Range("L132").Value = 22538257.2
$ws.Range("N132").Value = -22543317.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 48987.25
$ws.Range("J95").Value = 48987.25
$ws.Range("L95").Value = 48987.25
$ws.Range("N95").Value = -54479.25
$ws.Range("H102").Value = 3544.9092
$ws.Range("I102").Value = 1799
$ws.Range("K102").Value = 1799
$ws.Range("M102").Value = -177
$ws.Range("H113").Value = 508278.16
$ws.Range("I113").Value = 2002440
$ws.Range("K113").Value = 2002440
$ws.Range("M113").Value = -2000270
$ws.Range("H122").Value = 4455.222
$ws.Range("I122").Value = 2399
$ws.Range("K122").Value = 7197
$ws.Range("M122").Value = -4747
$ws.Range("H132").Value = 52118.477
$ws.Range("I132").Value = 5922.231
$ws.Range("J132").Value = 127187.375
$ws.Range("K132").Value = 17766.693
$ws.Range("L132").Value = 381562.125
$ws.Range("M132").Value = -15236.693
$ws.Range("N132").Value = -386622.125
$ws.Range("H136").Value = 54199.727
$ws.Range("J136").Value = 54199.727
$ws.Range("L136").Value = 162599.181
$ws.Range("N136").Value = -167699.181

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 1000
$ws.Range("M82").Value = -639
$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 1000
$ws.Range("M85").Value = 248
$ws.Range("H132").Value = 2879.4
$ws.Range("I132").Value = 2879.4
$ws.Range("K132").Value = 8638.200000000001
$ws.Range("M132").Value = -6108.200000000001
$ws.Range("H136").Value = 875679.0600000001
$ws.Range("I136").Value = 1181812.4
$ws.Range("J136").Value = 8301.166999999999
$ws.Range("K136").Value = 3545437.2
$ws.Range("L136").Value = 24903.501
$ws.Range("M136").Value = -3542887.2
$ws.Range("N136").Value = -30003.501

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 334266.34
$ws.Range("J96").Value = 2000
$ws.Range("L96").Value = 2000
$ws.Range("N96").Value = -4746
$ws.Range("H108").Value = 89403
$ws.Range("J108").Value = 89870.5
$ws.Range("L108").Value = 89870.5
$ws.Range("N108").Value = -97550.5
$ws.Range("H126").Value = 1699.2858
$ws.Range("I126").Value = 1699.2858
$ws.Range("K126").Value = 5097.857400000001
$ws.Range("M126").Value = -2627.857400000001
$ws.Range("H132").Value = 104790.5
$ws.Range("I132").Value = 4571.4287
$ws.Range("K132").Value = 13714.2861
$ws.Range("M132").Value = -11184.2861
$ws.Range("H138").Value = 250064860
$ws.Range("J138").Value = 250064860
$ws.Range("L138").Value = 250064860
$ws.Range("N138").Value = -250075140
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
